{"js": "// Apply the edits described in the commit: several wording / wording-expansion\n// changes throughout the \"Research Ethics Application Form\" document.\n// Each edit is performed as an independent search + replace so that runs\n// (and any surrounding run-level markup such as w:proofErr) outside the\n// matched text are left completely untouched.\n\nasync function replaceOnce(body, searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) \"How will participants be recruited...\" answer: drop \"game design\"\n//    before \"lecturers\" and add the \"during different guild sessions\" detail.\nawait replaceOnce(\n  body,\n  \"will ask the game design lecturers in the\",\n  \"will ask the lecturers in the\"\n);\n\nawait replaceOnce(\n  body,\n  \" academy for a random \",\n  \" academy during different guild sessions for a random \"\n);\n\n// 2) \"What will participants be asked to do?\" answer: replaced entirely with\n//    a description of the new website-based task. Use the paragraph object\n//    (rather than a body-wide search) so that the now-unused w:proofErr\n//    markers bracketing the old \"size  and\" run are removed together with\n//    the rest of the paragraph's old content.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\nconst levelsParagraph = paragraphs.items.find(\n  (p) =>\n    p.text.indexOf(\n      \"They will be required to design four levels, two levels of limited\"\n    ) !== -1\n);\nif (!levelsParagraph) {\n  throw new Error(\"Could not find the 'They will be required...' paragraph\");\n}\nlevelsParagraph.insertText(\n  \"They will be asked to access my secure HTTPS website and will be  ask to design five levels with the software on the website.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 3) \"What potential risks...\" (participants) answer: \"age\" -> \"date of birth\".\nawait replaceOnce(\n  body,\n  \"personal information such as age, gender etc. will not be needed.\",\n  \"personal information such as date of birth, gender etc. will not be needed.\"\n);\n\n// 4) Data storage answer: hard drive -> secure server, plus the new sentence\n//    about HTTPS encryption in transit.\nawait replaceOnce(\n  body,\n  \"The data will be stored on a hard drive, that only I will have access to.\",\n  \"The data will be stored straight onto a secure server that only I will have access to. This data will be anonymous and encrypted as it will be travelling via HTTPS GET and POST requests.\"\n);\n", "ps1": "# Apply the edits described in the commit: several wording / wording-expansion\n# changes throughout the \"Research Ethics Application Form\" document.\n# Each edit is an independent Find & Replace against the whole document\n# story so that text outside the matched range (including w:proofErr\n# markers around grammar-checked runs) is left completely untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\n# 1) \"How will participants be recruited...\" answer: drop \"game design\"\n#    before \"lecturers\" and add the \"during different guild sessions\" detail.\nReplace-Text \"will ask the game design lecturers in the\" \"will ask the lecturers in the\"\nReplace-Text \" academy for a random \" \" academy during different guild sessions for a random \"\n\n# 2) \"What will participants be asked to do?\" answer: replaced entirely with\n#    a description of the new website-based task.\nReplace-Text \"They will be required to design four levels, two levels of limited size  and two levels at the maximum size possible.\" \"They will be asked to access my secure HTTPS website and will be  ask to design five levels with the software on the website.\"\n\n# 3) \"What potential risks...\" (participants) answer: \"age\" -> \"date of birth\".\nReplace-Text \"personal information such as age, gender etc. will not be needed.\" \"personal information such as date of birth, gender etc. will not be needed.\"\n\n# 4) Data storage answer: hard drive -> secure server, plus the new sentence\n#    about HTTPS encryption in transit.\nReplace-Text \"The data will be stored on a hard drive, that only I will have access to.\" \"The data will be stored straight onto a secure server that only I will have access to. This data will be anonymous and encrypted as it will be travelling via HTTPS GET and POST requests.\"\n"}
